# Adds a new "2023" column (T) to the worksheet, mirroring the existing
# year columns (D:S, years 2007-2022). For every data row, column T gets
# the same number formatting/styling as column S (its immediate left
# neighbor) via a formats-only paste, and then the new 2023 value is
# written into it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number (3 = header/year row, 4-33 = data rows) -> new 2023 value
$newColumnValues = [ordered]@{
    3  = 2023
    4  = 52.734251206028382
    5  = 44.646801162600475
    6  = 60.998061560200554
    7  = 41.931627189714625
    8  = 38.177163051511151
    9  = 45.607453560981966
    10 = 50.172884880431361
    11 = 44.112367891063748
    12 = 56.155144351753421
    13 = 37.12775271808399
    14 = 26.579446704517768
    15 = 47.691579663423148
    16 = 49.966474107695483
    17 = 44.339536521432947
    18 = 55.468421253968863
    19 = 45.077411133103766
    20 = 42.341975649266388
    21 = 47.772457765110225
    22 = 46.127136558116561
    23 = 38.861148383596195
    24 = 53.259250196123595
    25 = 77.010722119335071
    26 = 59.189709472566221
    27 = 95.060584781341987
    28 = 61.976853895626128
    29 = 48.390234028455353
    30 = 78.407224173903401
    31 = 39.559787476030614
    32 = 41.059215262778977
    33 = 38.128175110356899
}

foreach ($row in $newColumnValues.Keys) {
    $sourceCell = $ws.Range("S$row")
    $targetCell = $ws.Range("T$row")

    # Copy column S's formatting (number format, font, alignment, ...)
    # onto column T before writing the value, so T looks like the other
    # year columns instead of picking up the sheet's default style.
    $sourceCell.Copy()
    $targetCell.PasteSpecial(-4122)  # xlPasteFormats

    $targetCell.Value = $newColumnValues[$row]
}
